$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 3
$ws.Range("G2").Value = 56.47185866666666
$ws.Range("H2").Value = 169.415576
$ws.Range("I2").Value = 0.8070274173741353
$ws.Range("J2").Value = 0.8070274173741354
$ws.Range("K2").Value = 3
$ws.Range("M2").Value = 12.35607266666667
$ws.Range("N2").Value = 37.068218
$ws.Range("Q2").Value = 697.7703893070632
$ws.Range("R2").Value = 6279.933503763567
$ws.Range("S2").Value = 0.8070274173741353
$ws.Range("T2").Value = 0.8070274173741354

$ws.Range("E3").Value = 3
$ws.Range("G3").Value = 10.23495333333333
$ws.Range("H3").Value = 30.70486
$ws.Range("I3").Value = 0.1462655586439962
$ws.Range("J3").Value = 0.1462655586439962
$ws.Range("K3").Value = 3
$ws.Range("M3").Value = 12.35607266666667
$ws.Range("N3").Value = 37.068218
$ws.Range("Q3").Value = 126.4638271266089
$ws.Range("R3").Value = 1138.17444413948
$ws.Range("S3").Value = 0.1462655586439962
$ws.Range("T3").Value = 0.1462655586439962

$ws.Range("E4").Value = 3
$ws.Range("G4").Value = 3.268330666666667
$ws.Range("H4").Value = 9.804992
$ws.Range("I4").Value = 0.04670702398186845
$ws.Range("J4").Value = 0.04670702398186846
$ws.Range("K4").Value = 3
$ws.Range("M4").Value = 12.35607266666667
$ws.Range("N4").Value = 37.068218
$ws.Range("Q4").Value = 40.38373121602845
$ws.Range("R4").Value = 363.453580944256
$ws.Range("S4").Value = 0.04670702398186845
$ws.Range("T4").Value = 0.04670702398186846
